# Generate Report for Handback
# Updates the Overview / zh-cn / de-de sheets of the localization-status
# workbook to reflect that both source files have been handed back and are
# now in sync with en-US: status text changes, "Latest Target File" /
# "Latest Handback File" hyperlink columns get populated, and the
# "Latest Handback DateTime" timestamps are recorded.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status text (column zh-cn / de-de) changes
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: 9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md
$wsZh.Range("B2").Value = $newStatus
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f960597d86a944c4693fc8fd8ac61f3cbe9e45a0/e2e/9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md",
    "",
    "",
    "9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a82679c1a26e11279d2bfeff6ac3178298eea6b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yufeih/9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.75fe3ad83d3b1f387cd7848e25355d8c209a2927.zh-cn.xlf",
    "",
    "",
    "9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.75fe3ad83d3b1f387cd7848e25355d8c209a2927.zh-cn.xlf"
) | Out-Null
$wsZh.Range("G2").Value = "2016-02-22 03:23:58"

# Row 3: b2b6df03-c520-49b2-90a1-848d14002757.md
$wsZh.Range("B3").Value = $newStatus
$wsZh.Hyperlinks.Add(
    $wsZh.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f960597d86a944c4693fc8fd8ac61f3cbe9e45a0/e2e/b2b6df03-c520-49b2-90a1-848d14002757.md",
    "",
    "",
    "b2b6df03-c520-49b2-90a1-848d14002757.md"
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a82679c1a26e11279d2bfeff6ac3178298eea6b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yufeih/b2b6df03-c520-49b2-90a1-848d14002757.7a03ac360e2d593e9b47d7a32f97c147aa08688c.zh-cn.xlf",
    "",
    "",
    "b2b6df03-c520-49b2-90a1-848d14002757.7a03ac360e2d593e9b47d7a32f97c147aa08688c.zh-cn.xlf"
) | Out-Null
$wsZh.Range("G3").Value = "2016-02-22 03:23:58"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: 9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md
$wsDe.Range("B2").Value = $newStatus
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f960597d86a944c4693fc8fd8ac61f3cbe9e45a0/e2e/9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md",
    "",
    "",
    "9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b51a16700aedf1c218d38f001c8ce77bcb8c02be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yufeih/9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.75fe3ad83d3b1f387cd7848e25355d8c209a2927.de-de.xlf",
    "",
    "",
    "9a22d94c-a763-4ece-9e4b-6e5c8efa1f55.75fe3ad83d3b1f387cd7848e25355d8c209a2927.de-de.xlf"
) | Out-Null
$wsDe.Range("G2").Value = "2016-02-22 03:24:20"

# Row 3: b2b6df03-c520-49b2-90a1-848d14002757.md
$wsDe.Range("B3").Value = $newStatus
$wsDe.Hyperlinks.Add(
    $wsDe.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f960597d86a944c4693fc8fd8ac61f3cbe9e45a0/e2e/b2b6df03-c520-49b2-90a1-848d14002757.md",
    "",
    "",
    "b2b6df03-c520-49b2-90a1-848d14002757.md"
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b51a16700aedf1c218d38f001c8ce77bcb8c02be/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yufeih/b2b6df03-c520-49b2-90a1-848d14002757.7a03ac360e2d593e9b47d7a32f97c147aa08688c.de-de.xlf",
    "",
    "",
    "b2b6df03-c520-49b2-90a1-848d14002757.7a03ac360e2d593e9b47d7a32f97c147aa08688c.de-de.xlf"
) | Out-Null
$wsDe.Range("G3").Value = "2016-02-22 03:24:20"
